$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add new header for the week of 28_01_2024
$ws.Range("F1").Value = "28_01_2024"

# Add the new week's sales figures for each recepcionista
$ws.Range("F2").Value = 838
$ws.Range("F3").Value = 720
$ws.Range("F4").Value = 1330
$ws.Range("F5").Value = 2851

# Update the selection to reflect where the user left off (H6)
$ws.Cells.Item(6, 8).Select() | Out-Null
